$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header rows (column A has the rotation title, style s=1) -- add "Price" header in column E
$headerRows = @(1, 9, 18, 27, 35, 43, 52, 61)
foreach ($r in $headerRows) {
    $ws.Cells.Item($r, 5).Value = "Price"
}

# Data rows -- add computed Price values in column E
$priceData = @{
    2 = 3300; 3 = 4125; 4 = 4125; 5 = 4125; 6 = 13200; 7 = 5000;
    10 = 4125; 11 = 3000; 12 = 6300; 13 = 5000; 14 = 6750; 15 = 5000; 16 = 3300;
    19 = 5000; 20 = 4125; 21 = 4125; 22 = 4125; 23 = 9000; 24 = 7875; 25 = 4250;
    28 = 4250; 29 = 3000; 30 = 10000; 31 = 4125; 32 = 5000; 33 = 7875;
    36 = 3000; 37 = 3000; 38 = 7875; 39 = 5000; 40 = 15000; 41 = 5000;
    44 = 3000; 45 = 4125; 46 = 7875; 47 = 5000; 48 = 15000; 49 = 5000; 50 = 2400;
    53 = 6250; 54 = 3000; 55 = 10000; 56 = 5000; 57 = 13200; 58 = 7875; 59 = 4125;
    62 = 6250; 63 = 4125; 64 = 7875; 65 = 4125; 66 = 12000; 67 = 7875;
}

foreach ($r in $priceData.Keys) {
    $ws.Cells.Item($r, 5).Value = $priceData[$r]
}

$ws.Range("E7").Select()
